# Applies the "Exit"/logout protocol additions described in the commit:
#   - documents the new "E" client->server exit command
#   - documents the matching server-side "Exit:" summary line

$d = $word.ActiveDocument

# --- Change 1: document the new "E" (exit) client command right after the
#     existing "g - get seats (...)" line, inside the "Client to Server
#     Protocols" paragraph. ---
$anchorA = $d.Content
$foundA = $anchorA.Find.Execute(
    "get seats (inside the vehicle information section).^l",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
if (-not $foundA) { throw "Could not locate the 'get seats' anchor text." }

$insA = $d.Range($anchorA.End, $anchorA.End)
$insA.InsertAfter("E " + [char]8211 + " exit from server log (protects the data of the server). " + [char]11)

# --- Change 2: document the matching "Exit:" entry at the end of the same
#     paragraph's protocol summary (after the "Get Seats Data:" line). ---
$anchorB = $d.Paragraphs.Item(9).Range
$foundB = $anchorB.Find.Execute(
    "view vehicle to get the seats)  g;len(vehicle_type);vehicle_type;len(vehicle_company);vehicle_company;len(vehicle_number); vehicle_number",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
if (-not $foundB) { throw "Could not locate the 'Get Seats Data' summary anchor text." }

$insB = $d.Range($anchorB.End, $anchorB.End)
$insB.InsertAfter("" + [char]11 + [char]11)

$insB.Collapse(0)
$insB.InsertAfter("Exit:")
$insB.Font.Underline = 1

$insB.Collapse(0)
$insB.InsertAfter(" E;len(username);username")
$insB.Font.Underline = 0
